$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "SumitDeshmukhIDNO1"
$ws.Range("B1").Value = "SumitDeshmukhPASSNO1"
$ws.Range("A2").Value = "SumitDeshmukhIDNO2"
$ws.Range("B2").Value = "SumitDeshmukhPASSNO2"
$ws.Range("A3").Value = "SumitDeshmukhIDNO3"
$ws.Range("B3").Value = "SumitDeshmukhPASSNO3"

$ws.Rows("4:6").Delete()

$ws.Range("A1:B3").Select()
